$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.011.38"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "2.461.87"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'486.94"
$ws.Range("E5").Value = "  +4.20%  "
$ws.Range("D6").Value = "'145.08"
$ws.Range("E6").Value = "  +9.28%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("D9").Value = "2.467.08"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").Value = "'5.81"
$ws.Range("E10").Value = "  +8.96%  "
$ws.Range("D11").Value = "'0.0968"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").Value = "'0.330"
$ws.Range("E12").Value = "  +4.31%  "
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").Value = "2.890.26"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "56.032.08"
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").Value = "'21.04"
$ws.Range("E16").Value = "  +6.36%  "
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "2.478.97"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "'4.50"
$ws.Range("E19").Value = "  +6.23%  "
$ws.Range("D20").Value = "'10.03"
$ws.Range("E20").Value = "  +4.55%  "
$ws.Range("D21").Value = "'316.32"
$ws.Range("E21").Value = "  +1.03%  "
$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "'5.78"
$ws.Range("E23").Value = "  +7.17%  "
$ws.Range("D24").Value = "'58.43"
$ws.Range("E24").Value = "  +2.91%  "
$ws.Range("D25").Value = "'0.410"
$ws.Range("E25").Value = "  +5.93%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'0.160"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").Value = "2.578.45"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "'7.67"
$ws.Range("E29").Value = "  +6.96%  "
$ws.Range("D30").Value = "0.0₃0781"
$ws.Range("E30").Value = "  +8.14%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'147.92"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").Value = "'18.20"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("E34").Value = "  +4.53%  "
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("E36").Value = "  +7.59%  "
$ws.Range("D37").Value = "'3.72"
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("E38").Value = "  +6.93%  "
$ws.Range("D39").Value = "'33.90"
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("D40").Value = "'3.51"
$ws.Range("E40").Value = "  +7.52%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'0.0553"
$ws.Range("E42").Value = "  +4.92%  "
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("E44").Value = "  +6.25%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'4.70"
$ws.Range("E45").Value = "  +11.18%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'259.88"
$ws.Range("E46").Value = "  +10.69%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0920"
$ws.Range("E47").Value = "  +3.97%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'10.20"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("E49").Value = "  +3.91%  "
$ws.Range("D50").Value = "'17.52"
$ws.Range("E50").Value = "  +4.92%  "
$ws.Range("D51").Value = "1.870.37"
$ws.Range("E51").Value = "  -3.84%  "
